# Add latest 816 (Versant Power: MPS + BHE subutilities) and CMP standard
# offer rate data for Jan 2022 - Dec 2022 (monthly, 1st-of-month dates),
# appended after the existing data (through row 361).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 12 month-start dates (Excel serials) used for every new block.
$dates = @(44562, 44593, 44621, 44652, 44682, 44713, 44743, 44774, 44805, 44835, 44866, 44896)

# ---------------------------------------------------------------------
# Block 1: rows 362-373 -> Versant Power / MPS, rate 0.110873
# Use the last existing Versant/MPS row (361) as the format template.
# ---------------------------------------------------------------------
$ws.Range("A361:D361").Copy()
$ws.Range("A362:D373").PasteSpecial(-4122)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(362 + $i, 1).Value = $dates[$i]
}
$ws.Range("B362:B373").Value = "Versant Power"
$ws.Range("C362:C373").Value = "MPS"
$ws.Range("D362:D373").Value = 0.110873

# ---------------------------------------------------------------------
# Block 2: rows 374-385 -> Versant Power / BHE, rate 0.11684
# Use the last existing Versant/BHE row (337) as the format template.
# ---------------------------------------------------------------------
$ws.Range("A337:D337").Copy()
$ws.Range("A374:D385").PasteSpecial(-4122)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(374 + $i, 1).Value = $dates[$i]
}
$ws.Range("B374:B385").Value = "Versant Power"
$ws.Range("C374:C385").Value = "BHE"
$ws.Range("D374:D385").Value = 0.11684

# ---------------------------------------------------------------------
# Block 3: rows 386-397 -> CMP (no subutility), rate 0.118161
# Use the last existing CMP row (313) as the format template; CMP rows
# never have a value in column C, so clear it back out after the paste.
# ---------------------------------------------------------------------
$ws.Range("A313:D313").Copy()
$ws.Range("A386:D397").PasteSpecial(-4122)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(386 + $i, 1).Value = $dates[$i]
}
$ws.Range("B386:B397").Value = "CMP"
$ws.Range("D386:D397").Value = 0.118161
$ws.Range("C386:C397").ClearContents()

# ---------------------------------------------------------------------
# Restore the selection/view state to reflect where editing finished.
# ---------------------------------------------------------------------
[void]$ws.Range("C400").Select()
